$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# YDS sheet: append new game's play-by-play yardage logs to the
# existing space-separated strings (R/P rows x OFF/DEF columns).
# -----------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")

$wsYDS.Range("B2").Value2 = $wsYDS.Range("B2").Value2 + " -3 4 5 5 8 5 1 5 0 0 1 1 9 3 3 2 4 0 2 5 17"
$wsYDS.Range("C2").Value2 = $wsYDS.Range("C2").Value2 + " 9 1 13 4 4 9 5 -1 6 -2 0 1 5 2 2 -1 7 1 5 11 3 9 26 -1 16 1 13 5 1 5 5 3 0 5 1 2 3 9"
$wsYDS.Range("B3").Value2 = $wsYDS.Range("B3").Value2 + " 10 7 18 20 9 12 8 6 0 5 14 8 6 24 5 6 2 11 10 38 9 10 9"
$wsYDS.Range("C3").Value2 = $wsYDS.Range("C3").Value2 + " 17 11 15 11 11 10 9 19 10 37 5 13 -1 18 5 4"

# -----------------------------------------------------------------
# OFF sheet: update season totals after the simulated Wild Card game.
# Row 2 = Home, Row 3 = Road
# -----------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")

$wsOFF.Range("C2").Value = 219
$wsOFF.Range("F2").Value = 61
$wsOFF.Range("G2").Value = 61
$wsOFF.Range("J2").Value = 36
$wsOFF.Range("L2").Value = 324
$wsOFF.Range("M2").Value = 212
$wsOFF.Range("O2").Value = 24
$wsOFF.Range("P2").Value = 13
$wsOFF.Range("Q2").Value = 592

$wsOFF.Range("B3").Value = 17
$wsOFF.Range("C3").Value = 212
$wsOFF.Range("E3").Value = 41
$wsOFF.Range("F3").Value = 143
$wsOFF.Range("G3").Value = 47
$wsOFF.Range("H3").Value = 27
$wsOFF.Range("I3").Value = 61
$wsOFF.Range("J3").Value = 64
$wsOFF.Range("N3").Value = 18

# -----------------------------------------------------------------
# DEF sheet: update season totals after the simulated Wild Card game.
# Row 2 = Home, Row 3 = Road
# -----------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")

$wsDEF.Range("B2").Value = 8
$wsDEF.Range("C2").Value = 210
$wsDEF.Range("D2").Value = 12
$wsDEF.Range("E2").Value = 16
$wsDEF.Range("F2").Value = 64
$wsDEF.Range("G2").Value = 41
$wsDEF.Range("I2").Value = 7
$wsDEF.Range("J2").Value = 27
$wsDEF.Range("L2").Value = 304
$wsDEF.Range("M2").Value = 182
$wsDEF.Range("Q2").Value = 582

$wsDEF.Range("C3").Value = 195
$wsDEF.Range("E3").Value = 40
$wsDEF.Range("F3").Value = 119
$wsDEF.Range("G3").Value = 38
$wsDEF.Range("H3").Value = 29
$wsDEF.Range("I3").Value = 72
$wsDEF.Range("J3").Value = 59

# -----------------------------------------------------------------
# ST sheet: update special-teams totals and append kick/punt
# distance logs for the new game.
# -----------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")

$wsST.Range("B2").Value = 101
$wsST.Range("D2").Value = 67
$wsST.Range("F2").Value = 354
$wsST.Range("G2").Value = 340
$wsST.Range("N2").Value = 69
$wsST.Range("O2").Value = 39

$wsST.Range("B3").Value = 84

$wsST.Range("B4").Value2 = $wsST.Range("B4").Value2 + " 55 60"
$wsST.Range("D3").Value2 = $wsST.Range("D3").Value2 + " 57 59 46 48 58"
$wsST.Range("D4").Value2 = $wsST.Range("D4").Value2 + " 11 0 0 0 0"
$wsST.Range("B5").Value2 = $wsST.Range("B5").Value2 + " 16 13"
$wsST.Range("D5").Value2 = $wsST.Range("D5").Value2 + " 0 0 5 0"
$wsST.Range("B6").Value2 = $wsST.Range("B6").Value2 + " 18 32"

# -----------------------------------------------------------------
# TURNS sheet: update turnover totals.
# -----------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")

$wsTURNS.Range("B2").Value = 4
$wsTURNS.Range("C2").Value = 12
$wsTURNS.Range("D2").Value = 12

$wsTURNS.Range("D3").Value = 10

# -----------------------------------------------------------------
# PEN sheet: update penalty totals.
# -----------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")

$wsPEN.Range("B2").Value = 19
$wsPEN.Range("D2").Value = 18

$wsPEN.Range("B3").Value = 32
$wsPEN.Range("D3").Value = 6
